$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Update the acquisition timestamp (column A) for all data rows (2-7)
# from the previous run's timestamp to the new run's timestamp.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-09-19 06:33:55"
}
